# Add the "2022-Q4" quarterly holdings sheet, inserted right after the
# "总计" (summary) sheet, and record its totals as a new row in "总计".

$wb = $excel.ActiveWorkbook
$sheetTotal = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert a new row at the top of the data in "总计" for 2022-Q4 and fill
#    it in (pushes 2022-Q3 .. 2021-Q2 down by one row).
# ---------------------------------------------------------------------------
$sheetTotal.Rows.Item(2).Insert()

# Re-use existing formatting instead of Excel's inherited-from-header guess.
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)
$sheetTotal.Range("B3:D3").Copy()
$sheetTotal.Range("B2:D2").PasteSpecial(-4122)

$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q4"
$sheetTotal.Range("C2").Value = 6
$sheetTotal.Range("D2").Value = 0.25

# ---------------------------------------------------------------------------
# 2. Create the new "2022-Q4" worksheet right after "总计" and populate it
#    with the quarter's fund-holdings detail (same layout as the other
#    quarter sheets).
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $sheetTotal)
$newSheet.Name = "2022-Q4"

# Grab an already-correctly-formatted quarter sheet to copy styles from
# (must be looked up by name -- any *positional* reference grabbed before
# the Add() above would now resolve to the freshly inserted sheet instead).
$refSheet = $wb.Worksheets.Item("2022-Q3")

$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$refSheet.Range("A2").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B..G hold numeric-looking text (fund codes with leading zeros,
# percentages kept at fixed precision, etc.) -- force Text format so Excel
# doesn't silently coerce them to numbers.
$newSheet.Range("B2:G7").NumberFormat = "@"
# ... except the two rows whose market-value is a genuine numeric 0.
$newSheet.Range("G6:G7").NumberFormat = "General"

$data = @(
  @("002450", "平安睿享文娱灵活配置混合A", "3.85", "94.03", "3.43", "0.1321", 9),
  @("002451", "平安睿享文娱灵活配置混合C", "1.97", "94.03", "3.43", "0.0676", 9),
  @("002307", "银华多元视野灵活配置混合", "1.52", "89.13", "2.09", "0.0318", 10),
  @("005251", "银华多元动力灵活配置混合", "0.43", "88.30", "3.29", "0.0141", 2),
  @("014745", "恒生前海兴享混合C", "0.00", "82.60", "3.87", 0, 10),
  @("014744", "恒生前海兴享混合A", "0.00", "82.60", "3.87", 0, 10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $rec = $data[$i]
  $newSheet.Range("A$row").Value = $i
  $newSheet.Range("B$row").Value = $rec[0]
  $newSheet.Range("C$row").Value = $rec[1]
  $newSheet.Range("D$row").Value = $rec[2]
  $newSheet.Range("E$row").Value = $rec[3]
  $newSheet.Range("F$row").Value = $rec[4]
  $newSheet.Range("G$row").Value = $rec[5]
  $newSheet.Range("H$row").Value = $rec[6]
}
